# ADD results from server
#
# The source data now includes an extra "gb" metric (inserted right after
# "eb") and an extra "btes" metric (inserted right before "ites"); the
# "gt"/"dgt" metrics are gone. The net column count on each year-sheet
# stays the same (15 columns, A:O). We rewrite the header row and the
# single data row (row 2) on every year-sheet with the new layout and
# the refreshed values pulled from the server.

$wb = $excel.ActiveWorkbook

$headers = @("eb","gb","hp","st","wi","ieh","chp","ac","ab_ct","ab_hp","cp_ct","cp_hp","ttes","btes","ites")

# One row of data per year-sheet, in worksheet order (2025, 2030, 2035, 2040, 2045, 2050).
$allData = @(
    @(31251.19287316165, 0, 386830.1019569611, 0, 2317792.09144148, 75465.07577201782, 0, 20274.22343033684, 0, 0, 0, 0, 0, 18910.44901441236, 15952.10825463477),
    @(24677.2233802937, 0, 1196479.899689134, 0, 0, 56155.0398571117, 0, 45281.94202016341, 0, 0, 0, 0, 0, 41319.33120982331, 34059.82347464918),
    @(193961.9401775642, 0, 754668.2160649784, 0, 0, 1140.740783696684, 0, 31499.18821270974, 0, 0, 0, 0, 0, 42646.53192577381, 24018.47504444163),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 9262.79801954542, 0),
    @(61364.47115330531, 1544.757918726473, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 22517.87701973088, 6650.12915920241),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
)

for ($sheetIndex = 1; $sheetIndex -le $wb.Worksheets.Count; $sheetIndex++) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    $row = $allData[$sheetIndex - 1]

    for ($i = 0; $i -lt $headers.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item(1, $col).Value = $headers[$i]
        $ws.Cells.Item(2, $col).Value = $row[$i]
    }
}
